$wb = $excel.ActiveWorkbook

# --- Sheet "2025": update row 2 values ---
$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("B2").Value = 0.01031349999999998
$ws2025.Range("E2").Value = 0.3889546774358236
$ws2025.Range("I2").Value = 0.3871912
$ws2025.Range("L2").Value = 0.5977993
$ws2025.Range("M2").Value = 0.08392083333333335
$ws2025.Range("N2").Value = 10.81225050560962
$ws2025.Range("O2").Value = 2.869420863597251

# --- Sheet "2030": update row 2 values ---
$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("A2").Value = 0.05319467292261965
$ws2030.Range("B2").Value = 0.08031067743582365
$ws2030.Range("E2").Value = 0.5551468225641764
$ws2030.Range("I2").Value = 0.9830738000000001
$ws2030.Range("L2").Value = 0.3190583500000002
$ws2030.Range("M2").Value = 0.1360119166666667
$ws2030.Range("N2").Value = 19.54848800605096
$ws2030.Range("O2").Value = 8.898376319884211
